$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativacao date (needs to stay plain text, not auto-converted to a date serial) ---
$scratch = $ws.Range("A1")
$origScratchStyle = $scratch.Style
$scratch.NumberFormat = "@"
$scratch.Value = '01/01/2021'
$scratch.Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").PasteSpecial(-4163) | Out-Null
$scratch.ClearContents()
$scratch.Style = $origScratchStyle
$scratch.Value = $null
$excel.CutCopyMode = 0

# --- Plain text replacements ---
$ws.Range("B10").Value = 'a) Fornecer ao aluno uma visão integrada do processo de desenvolvimento de produto, desde as etapas iniciais de geração da ideia, avaliação econômica e desenvolvimento do conceito do produto até a preparação da fábrica, produção e lançamento do produto. b) Apresentar como os principais conhecimentos da formação de engenheiro de produção podem ser aplicados no processo de desenvolvimento de produto. c) Propiciar uma experiência prática de projeto.'
$ws.Range("C10").Value = 'a) Fornecer ao aluno uma visão integrada do processo de desenvolvimento de produto, desde as etapas iniciais de geração da ideia, avaliação econômica e desenvolvimento do conceito do produto até a preparação da fábrica, produção e lançamento do produto. b) Apresentar como os principais conhecimentos da formação de engenheiro de produção podem ser aplicados no processo de desenvolvimento de produto. c) Propiciar uma experiência prática de projeto.'

$ws.Range("B13").Value = '5701460 - Antonio Iacono'
$ws.Range("C13").Value = '5701460 - Antonio Iacono'

$ws.Range("B14").Value = '1. Conceito de processo de negócio e modelo de referência; 2. Desenvolvimento de produto como um processo: Fases principais do processo de desenvolvimento de produto; 3. Definição e conceitos básicos de gerenciamento de projetos; 4. Planejamento estratégico do produto; 5. Planejamento do Projeto; 6. Projeto Informacional; 7. Projeto Conceitual; 8. Projeto Detalhado; 9. Lançamento e Pós-Desenvolvimento. 10. Empreendedorismo e processo de desenvolvimento de produto.'
$ws.Range("C14").Value = '1. Conceito de processo de negócio e modelo de referência; 2. Desenvolvimento de produto como um processo: Fases principais do processo de desenvolvimento de produto; 3. Definição e conceitos básicos de gerenciamento de projetos; 4. Planejamento estratégico do produto; 5. Planejamento do Projeto; 6. Projeto Informacional; 7. Projeto Conceitual; 8. Projeto Detalhado; 9. Lançamento e Pós-Desenvolvimento. 10. Empreendedorismo e processo de desenvolvimento de produto.'

$ws.Range("B15").Value = '1.Concept of business and process reference model; 2. Product development as a process: main phases of the product development process; 3. Definition and basic concepts of project management; 4. Strategic product planning; 5. Project planning; 6. information design; 7. conceptual design; 8. detailed design; 9. launch and post-development. 10. Entrepreneurship and product development process'
$ws.Range("C15").Value = '1.Concept of business and process reference model; 2. Product development as a process: main phases of the product development process; 3. Definition and basic concepts of project management; 4. Strategic product planning; 5. Project planning; 6. information design; 7. conceptual design; 8. detailed design; 9. launch and post-development. 10. Entrepreneurship and product development process'

$ws.Range("B16").Value = '1. Visão geral do processo de desenvolvimento de produto. Conceito de processo de negócio. Modelo de referência e desenvolvimento de produto como um processo. Fases principais do processo de desenvolvimento de produto. Definição e conceitos básicos de gerenciamento de projetos. 2. Apresentação de um modelo de referência para desenvolvimento de produto. Descrição do modelo de referência. Importância e uso do modelo de referência. Principais características do modelo de referência para desenvolvimento de produto. 3. Apresentação das fases do processo de desenvolvimento de produto. Nesta parte são apresentados os conceitos e metodologias principais de cada fase do desenvolvimento de produto. 3.1. Planejamento estratégico de produtos. Conteúdo sobre a importância e como preparar a estratégia de produto, realizar o gerenciamento da carteira de projetos e elaborar a Minuta do Projeto. 3.2. Planejamento do projeto. Conteúdo sobre técnicas de planejamento do projeto e análise de viabilidade econômica. 3.3. Projeto informacional. Conceitos gerais de pesquisa de mercado, tipos de fontes de dados para pesquisa de mercado, diferenças entre as fontes de dados, formas de utilização de pesquisas e aprofundamento no tipo de pesquisa baseado em enquete. 3.4. Projeto conceitual. Apresentam-se os diversos métodos criativos para concepção, tais como Diagrama de funções, Matriz Morfológica e Matriz de Pugh, e demonstra-se com exemplos como utilizá-los de forma a gerar e escolher uma concepção de produto robusta.3.5. Projeto detalhado e preparação da produção. Apresenta-se o objetivo e resumo das atividades realizadas nas duas fases. No pós-desenvolvimento discute-se especialmente os aspectos relacionados com o meio ambiente, como descarte e reuso de peças do produto. 4. Empreendedorismo e processo de desenvolvimento de produto. O ambiente de inovação e o processo de desenvolvimento de produtos: empresas nascentes de base tecnológica, startups e aceleradoras.'
$ws.Range("C16").Value = '1. Visão geral do processo de desenvolvimento de produto. Conceito de processo de negócio. Modelo de referência e desenvolvimento de produto como um processo. Fases principais do processo de desenvolvimento de produto. Definição e conceitos básicos de gerenciamento de projetos. 2. Apresentação de um modelo de referência para desenvolvimento de produto. Descrição do modelo de referência. Importância e uso do modelo de referência. Principais características do modelo de referência para desenvolvimento de produto. 3. Apresentação das fases do processo de desenvolvimento de produto. Nesta parte são apresentados os conceitos e metodologias principais de cada fase do desenvolvimento de produto. 3.1. Planejamento estratégico de produtos. Conteúdo sobre a importância e como preparar a estratégia de produto, realizar o gerenciamento da carteira de projetos e elaborar a Minuta do Projeto. 3.2. Planejamento do projeto. Conteúdo sobre técnicas de planejamento do projeto e análise de viabilidade econômica. 3.3. Projeto informacional. Conceitos gerais de pesquisa de mercado, tipos de fontes de dados para pesquisa de mercado, diferenças entre as fontes de dados, formas de utilização de pesquisas e aprofundamento no tipo de pesquisa baseado em enquete. 3.4. Projeto conceitual. Apresentam-se os diversos métodos criativos para concepção, tais como Diagrama de funções, Matriz Morfológica e Matriz de Pugh, e demonstra-se com exemplos como utilizá-los de forma a gerar e escolher uma concepção de produto robusta.3.5. Projeto detalhado e preparação da produção. Apresenta-se o objetivo e resumo das atividades realizadas nas duas fases. No pós-desenvolvimento discute-se especialmente os aspectos relacionados com o meio ambiente, como descarte e reuso de peças do produto. 4. Empreendedorismo e processo de desenvolvimento de produto. O ambiente de inovação e o processo de desenvolvimento de produtos: empresas nascentes de base tecnológica, startups e aceleradoras.'

$ws.Range("B17").Value = '1. overview of the process of product development. concept of business process. reference model and product development as a process. main phases of the product development process. definition and basic concepts of project management. 2. presentation of a reference model for product development. description of the reference model. importance and use of the reference model. main features of the reference model for product development. 3. presentation of the stages of the product development process. in this section are presented the main concepts and methodologies of each phase of the product development. 3.1. strategic planning of products. information about the importance and how to prepare the product strategy, carry out the project portfolio management and draw up the draft of the project. 3.2. project planning. information about project planning techniques and analysis of economic feasibility. 3.3. informational project. general concepts of market research, types of data sources for market research, differences between the data sources, forms of use of research and study on the type of survey-based research. 3.4. conceptual design. the various creative methods to design, such as functions, morphological matrix diagram and array of pugh, and demonstrates with examples how to use them to generate and choose a robust product design. 3.5. detailed design and preparation of production. the objective and summary of activities carried out in two phases. in the post-development discusses in particular the aspects relating to the environment, such as disposal and reuse of parts of the product. 4. entrepreneurship and product development process: innovation environment and product development process: technology-based enterprises, startups and business accelerators'
$ws.Range("C17").Value = '1. overview of the process of product development. concept of business process. reference model and product development as a process. main phases of the product development process. definition and basic concepts of project management. 2. presentation of a reference model for product development. description of the reference model. importance and use of the reference model. main features of the reference model for product development. 3. presentation of the stages of the product development process. in this section are presented the main concepts and methodologies of each phase of the product development. 3.1. strategic planning of products. information about the importance and how to prepare the product strategy, carry out the project portfolio management and draw up the draft of the project. 3.2. project planning. information about project planning techniques and analysis of economic feasibility. 3.3. informational project. general concepts of market research, types of data sources for market research, differences between the data sources, forms of use of research and study on the type of survey-based research. 3.4. conceptual design. the various creative methods to design, such as functions, morphological matrix diagram and array of pugh, and demonstrates with examples how to use them to generate and choose a robust product design. 3.5. detailed design and preparation of production. the objective and summary of activities carried out in two phases. in the post-development discusses in particular the aspects relating to the environment, such as disposal and reuse of parts of the product. 4. entrepreneurship and product development process: innovation environment and product development process: technology-based enterprises, startups and business accelerators'

$ws.Range("B19").Value = 'Provas e trabalho'
$ws.Range("C19").Value = 'Provas e trabalho'

$ws.Range("B20").Value = 'M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre desenvolvimento de produtosM = Média de aproveitamento do aluno'
$ws.Range("C20").Value = 'M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre desenvolvimento de produtosM = Média de aproveitamento do aluno'

$ws.Range("B21").Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
$ws.Range("C21").Value = 'MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.'
